$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "Checklist" -> "Session"
$ws.Name = "Session"

# Remove the two oldest log rows (old rows 2 & 3: Student IDs 231249 / 231999).
# Deleting row 2 twice shifts the remaining rows up by two.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Remove the trailing row (old row 8 / Student ID 235020), which now sits at row 6.
$ws.Rows.Item(6).Delete()

# Update the remaining four data rows (now rows 2-5): log Type changes from
# "Selection" to "Scan", and the User column switches from the admin email to
# the session hash.
$hash = "5edfa2692bdacc5e6ee805c626c50cb44cebb065f092d9a1067d89f74dacd326"
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 5).Value = "Scan"
    $ws.Cells.Item($r, 6).Value = $hash
}
